{"js": "// Replace the date line and every \"dividend\u00f7divisor=quotient, remainder\"\n// answer cell with its updated value. All old values are unique in the\n// document, so a literal (non-wildcard) body.search + replace per pair is\n// safe and order-independent.\nconst replacements = [\n  [\"2024-01-21 Sunday\", \"2024-01-22 Monday\"],\n  [\"251\u00f76=41, 5\", \"646\u00f73=215, 1\"],\n  [\"275\u00f72=137, 1\", \"554\u00f73=184, 2\"],\n  [\"348\u00f79=38, 6\", \"618\u00f73=206, 0\"],\n  [\"530\u00f76=88, 2\", \"539\u00f75=107, 4\"],\n  [\"680\u00f73=226, 2\", \"935\u00f79=103, 8\"],\n  [\"315\u00f79=35, 0\", \"640\u00f79=71, 1\"],\n  [\"781\u00f78=97, 5\", \"281\u00f75=56, 1\"],\n  [\"118\u00f79=13, 1\", \"318\u00f72=159, 0\"],\n  [\"449\u00f73=149, 2\", \"654\u00f79=72, 6\"],\n  [\"434\u00f75=86, 4\", \"289\u00f79=32, 1\"],\n  [\"231\u00f73=77, 0\", \"134\u00f77=19, 1\"],\n  [\"919\u00f72=459, 1\", \"613\u00f75=122, 3\"],\n  [\"981\u00f72=490, 1\", \"638\u00f79=70, 8\"],\n  [\"796\u00f78=99, 4\", \"221\u00f78=27, 5\"],\n  [\"353\u00f77=50, 3\", \"598\u00f79=66, 4\"],\n  [\"527\u00f76=87, 5\", \"788\u00f77=112, 4\"],\n  [\"838\u00f72=419, 0\", \"882\u00f72=441, 0\"],\n  [\"631\u00f73=210, 1\", \"587\u00f77=83, 6\"],\n  [\"773\u00f73=257, 2\", \"808\u00f73=269, 1\"],\n  [\"366\u00f78=45, 6\", \"238\u00f76=39, 4\"],\n  [\"197\u00f74=49, 1\", \"761\u00f73=253, 2\"],\n  [\"514\u00f72=257, 0\", \"699\u00f76=116, 3\"],\n  [\"869\u00f79=96, 5\", \"512\u00f74=128, 0\"],\n  [\"797\u00f74=199, 1\", \"683\u00f72=341, 1\"],\n  [\"611\u00f73=203, 2\", \"972\u00f74=243, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"dividend\u00f7divisor=quotient, remainder\"\n# answer cell with its updated value. All old values are unique in the\n# document, so a literal Find/Replace per pair (wrap = wdFindContinue) is\n# safe and order-independent.\n$replacements = @(\n    @('2024-01-21 Sunday', '2024-01-22 Monday'),\n    @('251\u00f76=41, 5',       '646\u00f73=215, 1'),\n    @('275\u00f72=137, 1',      '554\u00f73=184, 2'),\n    @('348\u00f79=38, 6',       '618\u00f73=206, 0'),\n    @('530\u00f76=88, 2',       '539\u00f75=107, 4'),\n    @('680\u00f73=226, 2',      '935\u00f79=103, 8'),\n    @('315\u00f79=35, 0',       '640\u00f79=71, 1'),\n    @('781\u00f78=97, 5',       '281\u00f75=56, 1'),\n    @('118\u00f79=13, 1',       '318\u00f72=159, 0'),\n    @('449\u00f73=149, 2',      '654\u00f79=72, 6'),\n    @('434\u00f75=86, 4',       '289\u00f79=32, 1'),\n    @('231\u00f73=77, 0',       '134\u00f77=19, 1'),\n    @('919\u00f72=459, 1',      '613\u00f75=122, 3'),\n    @('981\u00f72=490, 1',      '638\u00f79=70, 8'),\n    @('796\u00f78=99, 4',       '221\u00f78=27, 5'),\n    @('353\u00f77=50, 3',       '598\u00f79=66, 4'),\n    @('527\u00f76=87, 5',       '788\u00f77=112, 4'),\n    @('838\u00f72=419, 0',      '882\u00f72=441, 0'),\n    @('631\u00f73=210, 1',      '587\u00f77=83, 6'),\n    @('773\u00f73=257, 2',      '808\u00f73=269, 1'),\n    @('366\u00f78=45, 6',       '238\u00f76=39, 4'),\n    @('197\u00f74=49, 1',       '761\u00f73=253, 2'),\n    @('514\u00f72=257, 0',      '699\u00f76=116, 3'),\n    @('869\u00f79=96, 5',       '512\u00f74=128, 0'),\n    @('797\u00f74=199, 1',      '683\u00f72=341, 1'),\n    @('611\u00f73=203, 2',      '972\u00f74=243, 0')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
